$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B:F (rows 2-25)
$dataBF = @(
        @(1.02, 1.017559111064485, 1.023407764977146, 1.018939353641001, 1.015923187398002),
        @(1.02, 1.018641279331411, 1.02421460950286, 1.01985983518596, 1.017649398123509),
        @(1.02, 1.01934073110362, 1.024735482784014, 1.020455173958533, 1.018765352390039),
        @(1.02, 1.019634595154042, 1.024954169127289, 1.020705389772821, 1.019234264465888),
        @(1.02, 1.019683925405109, 1.024990870611201, 1.020747398351298, 1.01931298334031),
        @(1.02, 1.019344658455674, 1.024738406015585, 1.020458517607056, 1.018771618928847),
        @(1.02, 1.017924998380114, 1.023680692622323, 1.019250492143294, 1.016506784903975),
        @(1.02, 1.015417280839279, 1.02180757455596, 1.017119646514397, 1.012507646697246),
        @(1.02, 1.013741235065084, 1.020552532862007, 1.01569756431084, 1.009835468157138),
        @(1.02, 1.01301445309697, 1.02000758030485, 1.015081409581992, 1.00867681532703),
        @(1.02, 1.012744334376385, 1.01980493241845, 1.014852483023902, 1.008246191455799),
        @(1.02, 1.012802283023989, 1.01984841142601, 1.014901591255269, 1.00833857316546),
        @(1.02, 1.012992128270233, 1.019990834030468, 1.015062487659491, 1.008641224959),
        @(1.02, 1.013109076966806, 1.020078555016869, 1.015161613343776, 1.00882766544373),
        @(1.02, 1.013789446973555, 1.020588667627138, 1.0157384482495, 1.009912329867372),
        @(1.02, 1.01421594385352, 1.020908242356359, 1.016100177681017, 1.010592279975688),
        @(1.02, 1.0144646118186, 1.021094499161987, 1.016311131422917, 1.010988731485821),
        @(1.02, 1.014549384215175, 1.021157983236882, 1.016383054983534, 1.011123885773421),
        @(1.02, 1.014170195186135, 1.02087397009508, 1.016061371367202, 1.010519343555699),
        @(1.02, 1.012936228003843, 1.019948900433635, 1.015015109327488, 1.008552108500798),
        @(1.02, 1.01215945940398, 1.01936595126461, 1.014356939198574, 1.007313785808724),
        @(1.02, 1.01257132771519, 1.019675109287694, 1.014705880754276, 1.007970384649254),
        @(1.02, 1.014190867349998, 1.020889456705052, 1.016078906384106, 1.010552300852707),
        @(1.02, 1.016066322145164, 1.022292926547201, 1.017670783409176, 1.013542550847414)
    )

# New values for columns I:M (rows 2-25)
$dataIM = @(
        @(1.02707312312399, 1.022772538728435, 1.026239292359949, 1.02178408081615, 1.018776895882761),
        @(1.027276406361592, 1.023489949910739, 1.026853009402694, 1.022510164933683, 1.02030582566243),
        @(1.027405610608922, 1.023952868757091, 1.027248272090583, 1.022979096177411, 1.021293697580629),
        @(1.027459369549374, 1.024147171354864, 1.02741399762755, 1.023176022035219, 1.021708660131189),
        @(1.027468363163718, 1.024179777564817, 1.027441797725, 1.023209074304714, 1.021778314528212),
        @(1.027406331133, 1.023955466248429, 1.027250488262688, 1.022981728345308, 1.021299243652166),
        @(1.027142306924615, 1.023015259973391, 1.026447085118599, 1.022029650078783, 1.019293909375985),
        @(1.0266591824082, 1.02134852342094, 1.025017155716229, 1.020345070067533, 1.015748813466623),
        @(1.026325064581994, 1.02023057157659, 1.024054246803869, 1.019217304709761, 1.013377230324225),
        @(1.026177529600688, 1.019744853838168, 1.023635000857672, 1.018727834530254, 1.012348256564307),
        @(1.026122298464311, 1.019564188769725, 1.023478927405718, 1.018545850558692, 1.011965730410637),
        @(1.026134165199431, 1.019602953287805, 1.023512421407183, 1.018584894591283, 1.01204779818138),
        @(1.026172972953675, 1.019729925074897, 1.023622106860433, 1.018712795216108, 1.0123166434063),
        @(1.026196826689602, 1.019808123680505, 1.023689641689544, 1.018791576045142, 1.012482245199141),
        @(1.026334795704764, 1.020262772407365, 1.024082022208322, 1.019249765069332, 1.013445475649278),
        @(1.02642057408674, 1.020547521879609, 1.024327535457025, 1.019536868494491, 1.014049125688754),
        @(1.026470331271095, 1.020713453592687, 1.024470517350835, 1.019704221151301, 1.014401026184804),
        @(1.026487250401085, 1.020770005303728, 1.02451923289347, 1.019761265487823, 1.014520981803951),
        @(1.026411399418717, 1.020516987296423, 1.02430121716396, 1.019506076395546, 1.01398438034112),
        @(1.026161556913021, 1.019692541906217, 1.023589816813257, 1.018675136457506, 1.012237484085151),
        @(1.026001982505603, 1.019172745285161, 1.02314052380745, 1.018151690059562, 1.011137287829489),
        @(1.026086811947839, 1.019448435980173, 1.023378893243347, 1.018429274378744, 1.011720701573994),
        @(1.026415545908778, 1.020530785047862, 1.024313109953169, 1.019519990359228, 1.014013636588841),
        @(1.026786201861421, 1.021780605678944, 1.025388518674154, 1.020781399307345, 1.01666671213319)
    )

for ($i = 0; $i -lt $dataBF.Length; $i++) {
    $r = $i + 2
    $rowBF = $dataBF[$i]
    for ($j = 0; $j -lt $rowBF.Length; $j++) {
        $ws.Cells.Item($r, $j + 2).Value = $rowBF[$j]
    }
    $rowIM = $dataIM[$i]
    for ($j = 0; $j -lt $rowIM.Length; $j++) {
        $ws.Cells.Item($r, $j + 9).Value = $rowIM[$j]
    }
}

Write-Output "Updated vm_pu values for 380 kV case"
